$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 683.2308
$ws.Range("I53").Value = 297.57144
$ws.Range("J53").Value = 1133.1666
$ws.Range("K53").Value = 297.57144
$ws.Range("L53").Value = 1133.1666
$ws.Range("M53").Value = 339.42856
$ws.Range("N53").Value = -2407.1666
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = ""
$ws.Range("H98").Value = 4487.6665
$ws.Range("I98").Value = 4116.75
$ws.Range("J98").Value = 5058.3076
$ws.Range("K98").Value = 4116.75
$ws.Range("L98").Value = 5058.3076
$ws.Range("M98").Value = -2618.75
$ws.Range("N98").Value = -8054.3076
$ws.Range("H103").Value = 799.1429000000001
$ws.Range("J103").Value = 707.3077
$ws.Range("L103").Value = 2121.9231
$ws.Range("N103").Value = -3293.9231
$ws.Range("H122").Value = 4487.6665
$ws.Range("I122").Value = 4116.75
$ws.Range("J122").Value = 5058.3076
$ws.Range("K122").Value = 12350.25
$ws.Range("L122").Value = 15174.9228
$ws.Range("M122").Value = -9900.25
$ws.Range("N122").Value = -20074.9228
$ws.Range("H138").Value = 5714.535
$ws.Range("I138").Value = 1588.7894
$ws.Range("J138").Value = 7222.019
$ws.Range("K138").Value = 4766.3682
$ws.Range("L138").Value = 21666.057
$ws.Range("M138").Value = 373.6318000000001
$ws.Range("N138").Value = -31946.057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5523.381
$ws.Range("I32").Value = 5028.222
$ws.Range("K32").Value = 5028.222
$ws.Range("M32").Value = -4741.222
$ws.Range("H45").Value = 3082.8333
$ws.Range("I45").Value = 2099.4
$ws.Range("K45").Value = 2099.4
$ws.Range("M45").Value = -1722.4
$ws.Range("H63").Value = 1841.7142
$ws.Range("I63").Value = 1338.6
$ws.Range("J63").Value = 3099.5
$ws.Range("K63").Value = 1338.6
$ws.Range("L63").Value = 3099.5
$ws.Range("M63").Value = -652.5999999999999
$ws.Range("N63").Value = -4471.5
$ws.Range("H66").Value = 1841.7142
$ws.Range("I66").Value = 1338.6
$ws.Range("J66").Value = 3099.5
$ws.Range("K66").Value = 6693
$ws.Range("L66").Value = 15497.5
$ws.Range("M66").Value = -3261
$ws.Range("N66").Value = -22361.5
$ws.Range("H74").Value = 27779256
$ws.Range("I74").Value = 32258958
$ws.Range("K74").Value = 32258958
$ws.Range("M74").Value = -32258084
$ws.Range("H77").Value = 27779256
$ws.Range("I77").Value = 32258958
$ws.Range("K77").Value = 161294790
$ws.Range("M77").Value = -161290422
$ws.Range("H102").Value = 3413.6
$ws.Range("I102").Value = 3015.7693
$ws.Range("K102").Value = 3015.7693
$ws.Range("M102").Value = -1393.7693
$ws.Range("H110").Value = 8847.308000000001
$ws.Range("I110").Value = 8893.888999999999
$ws.Range("J110").Value = 8742.5
$ws.Range("K110").Value = 8893.888999999999
$ws.Range("L110").Value = 8742.5
$ws.Range("M110").Value = -6848.888999999999
$ws.Range("N110").Value = -12832.5
$ws.Range("H122").Value = 33336208
$ws.Range("I122").Value = 3041.875
$ws.Range("J122").Value = 166668860
$ws.Range("K122").Value = 9125.625
$ws.Range("L122").Value = 500006580
$ws.Range("M122").Value = -6675.625
$ws.Range("N122").Value = -500011480

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = ""
$ws.Range("H86").Value = 8196.846
$ws.Range("I86").Value = 11857.333
$ws.Range("J86").Value = 5059.2856
$ws.Range("K86").Value = 11857.333
$ws.Range("L86").Value = 5059.2856
$ws.Range("M86").Value = -10734.333
$ws.Range("N86").Value = -7305.2856
$ws.Range("H89").Value = 8196.846
$ws.Range("I89").Value = 11857.333
$ws.Range("J89").Value = 5059.2856
$ws.Range("K89").Value = 59286.665
$ws.Range("L89").Value = 25296.428
$ws.Range("M89").Value = -53670.665
$ws.Range("N89").Value = -36528.428
$ws.Range("H105").Value = 6851.9062
$ws.Range("I105").Value = 11930.272
$ws.Range("J105").Value = 4191.8096
$ws.Range("K105").Value = 11930.272
$ws.Range("L105").Value = 4191.8096
$ws.Range("M105").Value = -10183.272
$ws.Range("N105").Value = -7685.8096

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 11797.7
$ws.Range("I105").Value = 1933.25
$ws.Range("K105").Value = 1933.25
$ws.Range("M105").Value = -186.25
$ws.Range("H122").Value = 2772210.5
$ws.Range("I122").Value = 2184.1667
$ws.Range("K122").Value = 6552.500100000001
$ws.Range("M122").Value = -4102.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6191.4
$ws.Range("I70").Value = 4657.6665
$ws.Range("J70").Value = 6848.7144
$ws.Range("K70").Value = 4657.6665
$ws.Range("L70").Value = 6848.7144
$ws.Range("M70").Value = -4387.6665
$ws.Range("N70").Value = -7388.7144
$ws.Range("H73").Value = 6191.4
$ws.Range("I73").Value = 4657.6665
$ws.Range("J73").Value = 6848.7144
$ws.Range("K73").Value = 4657.6665
$ws.Range("L73").Value = 6848.7144
$ws.Range("M73").Value = -3721.6665
$ws.Range("N73").Value = -8720.714400000001
$ws.Range("H102").Value = 965
$ws.Range("I102").Value = 861.6
$ws.Range("K102").Value = 861.6
$ws.Range("M102").Value = 760.4
$ws.Range("I122").Value = 1887.5
$ws.Range("K122").Value = 5662.5
$ws.Range("M122").Value = -3212.5
$ws.Range("H126").Value = 51803.35
$ws.Range("I126").Value = 63986.062
$ws.Range("J126").Value = 3072.5
$ws.Range("K126").Value = 191958.186
$ws.Range("L126").Value = 9217.5
$ws.Range("M126").Value = -189488.186
$ws.Range("N126").Value = -14157.5
$ws.Range("H135").Value = 214999.75
$ws.Range("J135").Value = 214999.75
$ws.Range("L135").Value = 214999.75
$ws.Range("N135").Value = -225139.75
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8189.6
$ws.Range("I7").Value = 7987
$ws.Range("J7").Value = 9000
$ws.Range("K7").Value = 7987
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = -7875
$ws.Range("N7").Value = -9224
$ws.Range("H68").Value = 2370.9375
$ws.Range("I68").Value = 2332.9092
$ws.Range("J68").Value = 2454.6
$ws.Range("K68").Value = 2332.9092
$ws.Range("L68").Value = 2454.6
$ws.Range("M68").Value = -1583.9092
$ws.Range("N68").Value = -3952.6
$ws.Range("H71").Value = 2370.9375
$ws.Range("I71").Value = 2332.9092
$ws.Range("J71").Value = 2454.6
$ws.Range("K71").Value = 11664.546
$ws.Range("L71").Value = 12273
$ws.Range("M71").Value = -7920.546
$ws.Range("N71").Value = -19761
$ws.Range("H93").Value = 1518481.1
$ws.Range("I93").Value = 3666.5
$ws.Range("J93").Value = 3790703
$ws.Range("K93").Value = 3666.5
$ws.Range("L93").Value = 3790703
$ws.Range("M93").Value = -2418.5
$ws.Range("N93").Value = -3793199
$ws.Range("H126").Value = 8189.6
$ws.Range("I126").Value = 7987
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 23961
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -21491
$ws.Range("N126").Value = -31940
$ws.Range("H132").Value = 5234.8
$ws.Range("I132").Value = 2309.6155
$ws.Range("J132").Value = 24248.5
$ws.Range("K132").Value = 6928.8465
$ws.Range("L132").Value = 72745.5
$ws.Range("M132").Value = -4398.8465
$ws.Range("N132").Value = -77805.5
$ws.Range("H136").Value = 2508437.5
$ws.Range("I136").Value = 3336916.5
$ws.Range("J136").Value = 23000
$ws.Range("K136").Value = 10010749.5
$ws.Range("L136").Value = 69000
$ws.Range("M136").Value = -10008199.5
$ws.Range("N136").Value = -74100
$ws.Range("H139").Value = 59499.5
$ws.Range("I139").Value = 59499.5
$ws.Range("K139").Value = 59499.5
$ws.Range("M139").Value = -54359.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4768.909
$ws.Range("H62").Value = 3591.3333
$ws.Range("I62").Value = 3399.5
$ws.Range("J62").Value = 3975
$ws.Range("K62").Value = 3399.5
$ws.Range("L62").Value = 3975
$ws.Range("M62").Value = -2775.5
$ws.Range("N62").Value = -5223
$ws.Range("H65").Value = 3591.3333
$ws.Range("I65").Value = 3399.5
$ws.Range("J65").Value = 3975
$ws.Range("K65").Value = 16997.5
$ws.Range("L65").Value = 19875
$ws.Range("M65").Value = -13877.5
$ws.Range("N65").Value = -26115
$ws.Range("H132").Value = 4999.75
$ws.Range("I132").Value = 4999.75
$ws.Range("K132").Value = 14999.25
$ws.Range("M132").Value = -12469.25
